$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text (non-numeric-looking) cell updates: Coin name, Link, and combined label column
$ws.Range("E18").Value = '17OneONE'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("E43").Value = '42KickTokenKICK'
$ws.Range("E44").Value = '43LocalTradersLCTBestin24h'

# Price column (D) updates: force text format so numeric-looking strings are not
# converted to real numbers (source data stores these as text).
$priceCells = @{
    'D2' = '246.14'
    'D4' = '5.367'
    'D5' = '0.05733'
    'D6' = '6.474'
    'D7' = '3.139'
    'D9' = '0.8727'
    'D10' = '0.1379'
    'D11' = '0.06982'
    'D12' = '0.03135'
    'D13' = '0.02939'
    'D14' = '0.09417'
    'D15' = '3.744'
    'D16' = '0.001528'
    'D17' = '0.04705'
    'D18' = '0.0005978'
    'D19' = '0.006198'
    'D20' = '0.001238'
    'D21' = '0.004796'
    'D22' = '0.00008798'
    'D24' = '2.141'
    'D25' = '0.3175'
    'D26' = '0.1313'
    'D28' = '0.0002331'
    'D40' = '0.03718'
    'D41' = '0.1060'
    'D42' = '0.002390'
    'D43' = '0.003078'
    'D44' = '0.007509'
    'D45' = '0.00005279'
    'D47' = '0.3899'
    'D48' = '0.002604'
}
foreach ($ref in $priceCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$ref]
}
